$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 / Row 5: the two task names (and their "last updated" dates) were
# swapped -- "Task 3: Login Modal" now belongs with the earlier date and
# "Task 4: Register Modal" with the later date.
$ws.Range("A4").Value = "Task 3: Login Modal"
$ws.Range("A5").Value = "Task 4: Register Modal"

$ws.Range("D4").Value = (Get-Date -Year 2024 -Month 2 -Day 2).Date
$ws.Range("D5").Value = (Get-Date -Year 2024 -Month 5 -Day 2).Date

# Rows 8 / 9: "Create DB Tables" and "Create Factory for dummy data" are now
# finished -- mark them Complete (matching the style used by the other
# completed rows) and stamp a last-updated date (matching the date style
# already used in column D).
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B8").Value = "Complete"

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B9").Value = "Complete"

$ws.Range("D4").Copy()
$ws.Range("D8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D8").Value = (Get-Date -Year 2024 -Month 7 -Day 2).Date

$ws.Range("D4").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D9").Value = (Get-Date -Year 2024 -Month 7 -Day 2).Date

$excel.CutCopyMode = $false

# Restore the author's last selection.
$ws.Range("C12").Select()
